$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Update status of row 6 from "IN PROGRESS" to "DONE"
$ws.Range("F6").Value = "DONE"

# Update the selected cell on the sheet
$ws.Activate()
$ws.Range("F7").Select()
